# "fix(gui) step 1 and 2"
# Step 1: bump the printed/quoted date in A1 by one day.
# Step 2: update the two price cells (D35, D36) in the price table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: date ---
$ws.Range("A1").Value = 45309

# --- Step 2: prices ---
$ws.Range("D35").Value = 32935
$ws.Range("D36").Value = 7128

# The GUI re-touches every merged range on the sheet while re-rendering
# the layout after the edits above, which re-orders the worksheet's
# <mergeCells> list. Reproduce that by unmerging + re-merging each range
# in the resulting order.
$mergedRanges = @(
    "A10:D10",
    "B45:C45",
    "A11:D11",
    "A12:D12",
    "B37:C37",
    "A1:D1",
    "B35:C35",
    "B43:C43",
    "B44:C44",
    "B34:C34",
    "B41:C41",
    "B42:C42",
    "A9:D9",
    "B36:C36"
)

foreach ($ref in $mergedRanges) {
    $ws.Range($ref).UnMerge()
    $ws.Range($ref).Merge()
}
